# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G (header "K") held the old "Strike#" based values; these are
# regenerated here with the new "K" based values for each data row
# (rows 2-38, corresponding to A = 0..36).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(6,6,4,8,5,5,3,5,4,7,8,10,8,5,5,5,5,7,5,7,3,5,2,3,6,4,9,6,7,10,7,5,13,1,4,5,3)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("G$row").Value = $kValues[$i]
}
